$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the workbook default (unstyled) format,
# used to strip the quote-prefix style Excel applies when a numeric-
# looking string is forced to text below.
$defaultStyle = $ws.Range("A2").Style

$ws.Range("D2").Value = "29.146.15"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.822.81"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "'0.9982"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'234.26"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").Value = "'0.6018"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07054"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -5.02%  "
$ws.Range("D9").Value = "'0.2786"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").Value = "'23.42"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -6.08%  "
$ws.Range("D11").Value = "'0.07627"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "1.817.04"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "'4.788"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "'0.6274"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -6.93%  "
$ws.Range("D15").Value = "'0.000009912"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "2.064.50"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'78.45"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("D18").Value = "'5.841"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -6.10%  "
$ws.Range("D19").Value = "29.142.06"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "'225.69"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'1.0000"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'11.69"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").Value = "'6.963"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -5.17%  "
$ws.Range("D24").Value = "'0.9992"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'155.04"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").Value = "'7.982"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("D27").Value = "'0.1299"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "'0.06174"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -14.95%  "
$ws.Range("D31").Value = "'1.446"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "'3.828"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("D33").Value = "'3.789"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("D34").Value = "'1.120"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'1.735"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "'0.6361"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -8.41%  "
$ws.Range("D37").Value = "'2.538"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").Value = "1.211.87"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").Value = "'2.723"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").Value = "'0.01728"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("D41").Value = "'6.462"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -5.96%  "
$ws.Range("D42").Value = "'0.9059"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").Value = "'0.9992"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "1.978.88"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'100.32"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "'62.34"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -4.53%  "
$ws.Range("D47").Value = "'0.00000000115"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.588"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -6.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.460"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").Value = "'0.4552"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "'0.05501"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -2.84%  "
